$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the embedded pictures (Poza column used to hold images; now holds
# text asset paths instead) so the drawing layer is dropped.
$n = $ws.Shapes.Count
for ($i = $n; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Widen column E to fit the new "3D" asset-path content.
$ws.Columns.Item(5).ColumnWidth = 19.65

# Clear the old "Nr." numbering column (A2:A5).
$ws.Range("A2").Value = ""
$ws.Range("A3").Value = ""
$ws.Range("A4").Value = ""
$ws.Range("A5").Value = ""

# Poza column now stores the image asset file names as text.
$ws.Range("B2").Value = "assets/image1.jpg"
$ws.Range("B3").Value = "assets/image2.jpg"

# Updated dimensions text.
$ws.Range("C2").Value = "400 x 300 mm"
$ws.Range("C3").Value = "400 x 400 mm "

# 3D column now stores the 3D model asset file names as text.
$ws.Range("E2").Value = "assets/model1.glb"
$ws.Range("E3").Value = "assets/model2.glb"

# Update the remembered selection.
$null = $ws.Range("A5").Select()
